$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_pathway_genes")
$ws.Activate()

# Insert a new row at position 14, shifting existing rows 14-73 down to 15-74
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new gene entry
$ws.Range("A14").Value = "WP_085243324.1"
$ws.Range("B14").Value = "WP_085243324.1"
$ws.Range("C14").Value = "Chlorophyll"
$ws.Range("D14").Value = "chlI_bchI"
$ws.Range("E14").Value = "bchI"
$ws.Range("F14").Value = "S"
$ws.Range("G14").Value = 300
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = "Putative S-subunit in Actinobacteria [Mycobacterium europaeum]"

# Update the selection to match target state
$ws.Range("I2").Select()
